# Edit: L11 HW Assignment C
#  1. Table width: normalize "2500.0" pct -> "2500" pct (125%).
#  2. Make the first table row a repeating header row (w:tblHeader).
#  3. Trim the "Sketch the t-distribution using the t-distribution applet."
#     clause from problem statement in the numbered list (question about
#     finding the P-value).

$d = $word.ActiveDocument

# --- 1 & 2: table formatting -------------------------------------------------
$t = $d.Tables.Item(1)

# Re-assert the table's preferred width as 125% (2500 in fiftieths-of-a-
# percent units) so it serializes as an integer percentage.
$t.PreferredWidthType = 2
$t.PreferredWidth = 125

# Flag the first row as a heading row that repeats on each page.
$t.Rows.Item(1).HeadingFormat = $true

# --- 3: text trim -------------------------------------------------------------
$d.Content.Find.Execute(
    "Find the P-value and compare it to the level of significance. Sketch the t-distribution using the t-distribution applet.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Find the P-value and compare it to the level of significance.",
    2
)
